$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the diff: crypto price/volume table refresh (GitHub Actions run).
# Cells whose new value would be misread as a number by Excel's type-inference
# (plain decimals like "1.00", "0.999", "0.626", etc.) are forced to Text format
# first so the exact original string (incl. trailing zeros) round-trips unchanged.

$ws.Range("D2").Value = '63.637.53'
$ws.Range("E2").Value = '  -2.92%  '
$ws.Range("D3").Value = '2.604.50'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.42'
$ws.Range("E5").Value = '  -4.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.34'
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.628'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -5.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.78'
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.14'
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("D14").Value = '3.077.02'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("E15").Value = '  -7.68%  '
$ws.Range("D16").Value = '63.477.89'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").Value = '2.577.94'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.04'
$ws.Range("E18").Value = '  -4.00%  '
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.29'
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.93'
$ws.Range("E23").Value = '  -2.71%  '
$ws.Range("E24").Value = '  +2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000108'
$ws.Range("E25").Value = '  -4.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.24'
$ws.Range("E26").Value = '  -3.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '577.21'
$ws.Range("E27").Value = '  +9.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.57'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("E33").Value = '  -3.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.45'
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("E36").Value = '  -2.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.79'
$ws.Range("E37").Value = '  -2.63%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.39'
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.85'
$ws.Range("E40").Value = '  -3.75%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.35'
$ws.Range("E42").Value = '  -2.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '155.85'
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("E44").Value = '  +4.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.94'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.87'
$ws.Range("E46").Value = '  +1.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0591'
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.101'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.626'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.91'
$ws.Range("E51").Value = '  -3.95%  '
